$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.65198540687561
$ws.Range("B1").Value = 3.472146272659302
$ws.Range("C1").Value = 4.003684997558594
$ws.Range("D1").Value = 1.267269611358643
$ws.Range("E1").Value = 0.7433198690414429
